$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: new control entry at "poli" (Saint-Avold SA24 line), driver Bangoura
$ws.Range("A4").Value = "19/01/2026"
$ws.Range("B4").Value = "13:55"
$ws.Range("C4").Value = "13:55"
$ws.Range("D4").Value = "poli"
$ws.Range("E4").Value = "Bangoura"
$ws.Range("I4").Value = "Non observable"
$ws.Range("K4").Value = "tufd"
$ws.Range("L4").Value = "rgeFluo57"
$ws.Range("N4").Value = "Sa"
$ws.Range("T4").Value = "SA24"
$ws.Range("X4").Value = "beau"
$ws.Range("Y4").Value = 10235
$ws.Range("Z4").Value = "Conforme"
$ws.Range("AA4").Value = "Conforme"
$ws.Range("AB4").Value = "Conforme"
$ws.Range("AC4").Value = "Conforme"
$ws.Range("AD4").Value = "Conforme"
$ws.Range("AE4").Value = "Conforme"
$ws.Range("AF4").Value = "Propre"
$ws.Range("AG4").Value = "ras"
$ws.Range("AH4").Value = "Conforme"
$ws.Range("AI4").Value = "Conforme"
$ws.Range("AJ4").Value = "Conforme"
$ws.Range("AK4").Value = "Propre"
$ws.Range("AL4").Value = "Propre"
$ws.Range("AM4").Value = "Propre"
$ws.Range("AN4").Value = "Propre"
$ws.Range("AO4").Value = "ras"
$ws.Range("AP4").Value = 4
$ws.Range("AQ4").Value = 0
$ws.Range("AR4").Value = "BANGOURA"

# Row 5: new control entry at "poli" (Abris bus), driver Bangoura
$ws.Range("A5").Value = "19/01/2026"
$ws.Range("B5").Value = "14:07"
$ws.Range("C5").Value = "14:05"
$ws.Range("D5").Value = "poli"
$ws.Range("E5").Value = "Bangoura"
$ws.Range("F5").Value = "Conforme"
$ws.Range("G5").Value = "Conforme"
$ws.Range("H5").Value = "Conforme"
$ws.Range("I5").Value = "Abris bus"
$ws.Range("J5").Value = "Conforme"
$ws.Range("K5").Value = "RAS"
$ws.Range("L5").Value = "forbus"
$ws.Range("X5").Value = "beau"
$ws.Range("Y5").Value = 100125
$ws.Range("Z5").Value = "Conforme"
$ws.Range("AA5").Value = "Conforme"
$ws.Range("AB5").Value = "Conforme"
$ws.Range("AC5").Value = "Conforme"
$ws.Range("AD5").Value = "Conforme"
$ws.Range("AE5").Value = "Conforme"
$ws.Range("AF5").Value = "Propre"
$ws.Range("AG5").Value = "RAS"
$ws.Range("AH5").Value = "Conforme"
$ws.Range("AI5").Value = "Conforme"
$ws.Range("AJ5").Value = "Conforme"
$ws.Range("AK5").Value = "Propre"
$ws.Range("AL5").Value = "Propre"
$ws.Range("AM5").Value = "Propre"
$ws.Range("AN5").Value = "Propre"
$ws.Range("AO5").Value = "RAS"
$ws.Range("AP5").Value = 30
$ws.Range("AQ5").Value = 1
$ws.Range("AR5").Value = "BANGOURA"
